# Sprint backlog fix: the "Story Type" and "Story Points" columns (C and D)
# had been swapped for every data row, and the Story Type value "M" should
# actually read "F,T" (the unused "H" value is no longer needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-7 hold the sprint backlog entries.
for ($row = 2; $row -le 7; $row++) {
    # Capture the current (mis-placed) Story Points numeric value, which
    # currently lives in column C.
    $points = $ws.Cells.Item($row, 3).Value2

    # Column C becomes the corrected Story Type text, column D becomes the
    # numeric Story Points value that used to live in C.
    $ws.Cells.Item($row, 3).Value = "F,T"
    $ws.Cells.Item($row, 4).Value = $points
}

# Reflect the author's final cell selection.
$ws.Range("C10").Select()
